$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Update row 78 (D78:M78) with new indicator text ---
$ws.Range('D78').Value = 'Ausgebildetes Personal zur Kleinwaffenkontrolle'
$ws.Range('E78').Value = 'XXXAusgebildetes Personal zur Kleinwaffenkontrolle'
$ws.Range('F78').Value = 'Von der Bundeswehr ausgebildetes Personal zur Stärkung der Kleinwaffenkontrolle und Munitionssicherheit'
$ws.Range('G78').Value = 'XXXVon der Bundeswehr ausgebildetes Personal zur Stärkung der Kleinwaffenkontrolle und Munitionssicherheit'
$ws.Range('H78').Value = 'Von 2025 bis 2030 Ausbildung von mindestens 1.000 Personen durch Expertinnen und Experten der Bundeswehr'
$ws.Range('I78').Value = 'XXXVon 2025 bis 2030 Ausbildung von mindestens 1.000 Personen durch Expertinnen und Experten der Bundeswehr'
$ws.Range('J78').Value = 'mindestens 1.000 Personen von 2025 bis 2030'
$ws.Range('K78').Value = 'XXXmindestens 1.000 Personen von 2025 bis 2030'
$ws.Range('L78').Value = 'Ausgebildetes Personal zur Kleinwaffenkontrolle'
$ws.Range('M78').Value = 'XXXAusgebildetes Personal zur Kleinwaffenkontrolle'

# --- Step 2: Shift rows 81-83 down to 82-84 (values only; style already s=4 on existing rows 82,83) ---
# Row 83 -> Row 84 (row 84 doesn't exist yet, needs format copy)
$ws.Range('A84').Value = 'Z17_B03_P01_Ib01_I01'
$ws.Range('B84').Value = 'Z17_B03_P01_Ib01'
$ws.Range('C84').Value = '17.3'
$ws.Range('D84').Value = 'Einfuhren aus LDCs'
$ws.Range('E84').Value = 'Imports coming from LDCs'
$ws.Range('F84').Value = 'Einfuhren aus am wenigsten entwickelten Ländern'
$ws.Range('G84').Value = 'Imports from least developed countries'
$ws.Range('H84').Value = 'Steigerung des Anteils um 100 % bis 2030 gegenüber 2014'
$ws.Range('I84').Value = 'Increase the proportion by 100 % by 2030, compared to 2014'
$ws.Range('J84').Value = 'Steigerung des Anteils um 100 % bis 2030 gegenüber 2014'
$ws.Range('K84').Value = 'increase by 100 % by 2030 compared to 2014'
$ws.Range('L84').Value = 'Einfuhren aus am wenigsten entwickelten Ländern'
$ws.Range('M84').Value = 'Imports from least developed countries'

# Copy formatting from row 83 into new row 84 (row 84 has no prior style)
$ws.Range('A83:M83').Copy()
$ws.Range('A84:M84').PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 82 -> Row 83 (row 83 already has style s=4, just overwrite values)
$ws.Range('A83').Value = 'Z17_B02_P01_Ib01_I01'
$ws.Range('B83').Value = 'Z17_B02_P01_Ib01'
$ws.Range('C83').Value = '17.2'
$ws.Range('D83').Value = 'Studierende und Forschende aus Entwicklungslenländern und LDCs'
$ws.Range('E83').Value = 'Students and researchers from developing countries and LDCs'
$ws.Range('F83').Value = 'Anzahl der Studierenden und Forschenden aus Entwicklungsländern sowie aus am wenigsten entwickelten Ländern pro Jahr'
$ws.Range('G83').Value = 'Number of students and researchers from developing countries and least developed countries per year'
$ws.Range('H83').Value = 'Steigerung der Anzahl um 10 % von 2015 bis 2020, anschließend Verstetigung'
$ws.Range('I83').Value = 'Increase the number by 10 % from 2015 to 2020, then stabilised'
$ws.Range('J83').Value = 'Steigerung um 10 % von 2015 bis 2020, anschließend Verstetigung'
$ws.Range('K83').Value = 'increase by 10 % from 2015 to 2020, then stabilised'
$ws.Range('L83').Value = 'Anzahl der Studierenden und Forschenden aus Entwicklungsländern sowie aus am wenigsten entwickelten Ländern pro Jahr'
$ws.Range('M83').Value = 'Number of students and researchers from developing countries and least developed countries per year'

# Row 81 (old) -> Row 82 (row 82 already has style s=4, just overwrite values)
$ws.Range('A82').Value = 'Z17_B01_P01_Ib01_I01'
$ws.Range('B82').Value = 'Z17_B01_P01_Ib01'
$ws.Range('C82').Value = '17.1'
$ws.Range('D82').Value = 'Entwicklungsausgaben'
$ws.Range('E82').Value = 'Development assistance'
$ws.Range('F82').Value = 'Anteil öffentlicher Entwicklungsausgaben am Bruttonationaleinkommen'
$ws.Range('G82').Value = 'Official development assistance as a proportion of gross national income'
$ws.Range('H82').Value = 'Steigerung des Anteils auf 0,7 % des  Bruttonationaleinkommens bis 2030'
$ws.Range('I82').Value = 'Increase the proportion to 0.7 % of gross national income by 2030'
$ws.Range('J82').Value = 'Steigerung des Anteils auf 0,7 % des BNE bis 2030'
$ws.Range('K82').Value = 'increase to 0.7 % of GNI by 2030'
$ws.Range('L82').Value = 'Anteil öffentlicher Entwicklungsausgaben am Bruttonationaleinkommen'
$ws.Range('M82').Value = 'Official development assistance as a proportion of gross national income'

# --- Step 3: Write new FATF row into row 81 (row 81 already has style s=4, just overwrite values) ---
$ws.Range('A81').Value = 'Z16_B04_P01_IB01_I01'
$ws.Range('B81').Value = 'Z16_B04_P01_IB01'
$ws.Range('C81').Value = '16.4'
$ws.Range('D81').Value = 'Financial Action Task Force (FATF) rating effectiveness'
$ws.Range('E81').Value = 'Financial Action Task Force (FATF) rating effectiveness'
$ws.Range('F81').Value = 'Financial Action Task Force (FATF) rating effectiveness'
$ws.Range('G81').Value = 'Financial Action Task Force (FATF) rating effectiveness'
$ws.Range('H81').Value = 'Verbesserung auf mindestens 8 von 11 möglichen Punkten bis 2029'
$ws.Range('I81').Value = 'XXXVerbesserung auf mindestens 8 von 11 möglichen Punkten bis 2029'
$ws.Range('J81').Value = 'Verbesserung auf mindestens 8 von 11 möglichen Punkten bis 2029'
$ws.Range('K81').Value = 'XXXVerbesserung auf mindestens 8 von 11 möglichen Punkten bis 2029'
$ws.Range('L81').Value = 'Financial Action Task Force (FATF) rating effectiveness'
$ws.Range('M81').Value = 'Financial Action Task Force (FATF) rating effectiveness'

# --- Step 4: Update column H width ---
# Target stored width is 39.6953125 (exact OOXML value). The COM ColumnWidth
# setter in this runtime quantizes to whole-pixel units (multiples of 1/7),
# so the exact value is not representable; 39 lands on the closest
# reachable width (39.714285714285715 => 278px @ MDW=7).
$ws.Columns('H:H').ColumnWidth = 39

Write-Output 'Done'
